# Nacimientos_por_edad_CR.xlsx
# "Cambio a frecuencia relativa de los datos"
# Convert the absolute birth counts (by mother's age group, 2003 vs 2023)
# into relative frequencies (percentages, comma as decimal separator).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row by row: 2003 column (B), then 2023 column (C) ---------------
$ws.Range("B2").Value = "20,5"
$ws.Range("C2").Value = "9,1"
$ws.Range("B3").Value = "30,3"
$ws.Range("C3").Value = "23,7"
$ws.Range("B4").Value = "23,2"
$ws.Range("C4").Value = "27,1"
$ws.Range("B5").Value = "15,8"
$ws.Range("C5").Value = "22,9"
$ws.Range("B6").Value = "8,1"
$ws.Range("C6").Value = "13,6"
$ws.Range("B7").Value = "2,0"
$ws.Range("C7").Value = "3,5"
$ws.Range("B8").Value = "0,1"
$ws.Range("C8").Value = "0,2"

# Column widths now that the data is narrower percentage strings.
$ws.Columns("A").ColumnWidth = 16
$ws.Columns("B:C").ColumnWidth = 4.166666666666667

# Move the active selection like the author left it.
$ws.Range("G4").Select()
